$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Compañias")

# Update the placeholder text used for the price-list column (E4)
$ws.Range("E4").Value = "{{item.PrecioLista}}"

# Move the active selection to E4 (single cell, matching saved view state)
$ws.Range("E4").Select()
